# loop-multiple-rows-and-nested-multiple-loops.xlsx
#
# Adds a new "Collaborators:" label row plus a new dynamic "cells" loop row
# right after the existing categories/tags block (rows 2-4), pushing every
# row below it down by two rows. Also normalizes the (now unused) bold-ish
# "applyFont" style that used to sit on the `{{#each reviews}}` and
# `{{#each stars}}...` rows back to the default style, and moves the
# worksheet selection to the newly inserted loop cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 5:6 (everything at/after row 5 shifts down by 2)
$ws.Range("A5:A6").EntireRow.Insert()

# Fill the new rows - set B6 first so the new loop-tag string is added to the
# shared string table before the "Collaborators:" label (matches save order).
$ws.Range("B6").Value = "{{#each cells=collaborators}}{{this}}{{/each}}"
$ws.Range("B5").Value = "Collaborators:"
$ws.Range("B5").Font.Bold = $true

# These two rows (originally the `{{#each reviews}}` and the
# `{{#each stars}}{{value}} - {{time}}` rows, now shifted to 14 and 18) had a
# stray "applyFont" style applied with no visible effect - reset to Normal.
$ws.Range("B14").Style = "Normal"
$ws.Range("B18").Style = "Normal"

# Move the selection to the newly added loop cell
$ws.Range("B6").Select() | Out-Null
